$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collect all cells whose new "Price" text would otherwise be auto-converted
# to a floating point number by Excel, and force them to remain plain text
# so the stored value exactly matches the original literal string.
$textCells = @("D4", "D5", "D6", "D7", "D10", "D11", "D12", "D15", "D17", "D19", "D21", "D22", "D23", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D41", "D42", "D44", "D47", "D49", "D51")
foreach ($c in $textCells) {
  $ws.Range($c).NumberFormat = "@"
}

# --- Row 37 / 38: Celestia and LidoDAOToken swapped position, with updated values ---
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "18.19"
$ws.Range("E37").Value = "  -4.35%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "3.16"
$ws.Range("E38").Value = "  -7.54%  "

# --- Price and Volume(1h) updates for all other rows ---
$ws.Range("D2").Value = "42.430.94"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.528.97"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "309.46"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").Value = "99.25"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Value = "35.79"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "0.0803"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "7.33"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "2.914.72"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "15.66"
$ws.Range("E15").Value = "  +3.97%  "
$ws.Range("D16").Value = "2.506.66"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "0.820"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "42.412.26"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").Value = "12.19"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").Value = "69.12"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "243.85"
$ws.Range("E23").Value = "  -3.62%  "
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "25.87"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").Value = "2.33"
$ws.Range("E28").Value = "  -4.64%  "
$ws.Range("D29").Value = "39.24"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").Value = "10.11"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").Value = "156.06"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "5.72"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("E33").Value = "  +14.41%  "
$ws.Range("D34").Value = "0.0793"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("E35").Value = "  -3.17%  "
$ws.Range("D36").Value = "2.01"
$ws.Range("E36").Value = "  -5.20%  "
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "4.31"
$ws.Range("E41").Value = "  +10.07%  "
$ws.Range("D42").Value = "21.68"
$ws.Range("E42").Value = "  -3.95%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "3.31"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").Value = "1.962.38"
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").Value = "8.91"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").Value = "2.772.03"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "81.05"
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").Value = "0.854"
$ws.Range("E51").Value = "  +9.77%  "
